$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: column D (Price) and column E (Volume 1h)
# text values. Use a leading apostrophe to force text interpretation (matches the
# original inlineStr cell type), then reset the style so no numFmt/quotePrefix
# style gets attached to the cell (keeps cells style-less, as in the source).

$ws.Range("D2").Value = "'29.202.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.855.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'241.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.6971"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.07756"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "'23.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "'0.07805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").Value = "'1.858.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'5.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'91.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "'0.6861"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'6.516"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "'0.000008444"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "'29.207.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'248.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'2.106.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'7.511"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'0.1493"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").Value = "'161.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "'8.856"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "'1.558"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'1.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'0.05205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.01863"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "'1.224.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "'2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'0.8995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'109.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'5.519"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.11%  "
$ws.Range("D45").Value = "'2.005.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").Value = "'65.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.41%  "
$ws.Range("D48").Value = "'0.5179"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'9.514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").Value = "'7.033"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
